$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new/renamed column labels ---
$ws.Range("C1").Value = "REER=Mean (entire period)"
$ws.Range("D1").Value = "REER=Mean (before EAEU)"
$ws.Range("E1").Value = "REER=Mean (after EAEU)"
$ws.Range("F1").Value = "Mean Reversion (before EAEU)"
$ws.Range("G1").Value = "Mean Reversion (after EAEU)"
$ws.Range("H1").Value = "Mean Reversion (entire period)"

# --- Columns C and D: swap the "entire period" / "before EAEU" mean values ---
$ws.Range("C2").Value = 0.3399999999999999
$ws.Range("D2").Value = 0.3199999999999999
$ws.Range("C3").Value = 0.4599999999999999
$ws.Range("D3").Value = 0.5699999999999997
$ws.Range("C4").Value = 0.5299999999999998
$ws.Range("D4").Value = 0.5799999999999997
$ws.Range("C5").Value = 0.4199999999999998
$ws.Range("D5").Value = 0.3899999999999999
$ws.Range("C6").Value = 0.4099999999999998
$ws.Range("D6").Value = 0.3599999999999999
$ws.Range("C7").Value = 0.4699999999999999
$ws.Range("D7").Value = 0.4799999999999998
$ws.Range("C8").Value = 0.4699999999999999
$ws.Range("D8").Value = 0.6799999999999997
$ws.Range("C9").Value = 0.4099999999999998
$ws.Range("D9").Value = 0.6199999999999998

# --- Column E: now holds the new "after EAEU" mean values (was text bucket before) ---
$ws.Range("E2").Value = 0.7299999999999996
$ws.Range("E3").Value = 0.1899999999999999
$ws.Range("E4").Value = 0.2499999999999999
$ws.Range("E5").Value = 0.6599999999999997
$ws.Range("E6").Value = 0.8199999999999996
$ws.Range("E7").Value = 0.4099999999999998
$ws.Range("E8").Value = 0.18
$ws.Range("E9").Value = 0.2

# --- Column F: write in row order so new shared strings are appended in the right sequence ---
$ws.Range("F2").Value = "N/A"
$ws.Range("F3").Value = "0.5-0.9*"
$ws.Range("F4").Value = "0.2"
$ws.Range("F5").Value = "0.1"
$ws.Range("F6").Value = "0.6, 0.8"
$ws.Range("F7").Value = "N/A"
$ws.Range("F8").Value = "0.1-0.4*"
$ws.Range("F9").Value = "0.4-0.7*"

# --- Column G (new column) ---
$ws.Range("G2").Value = "N/A"
$ws.Range("G3").Value = "N/A"
$ws.Range("G4").Value = "0.1-0.4*"
$ws.Range("G5").Value = "0.7,0.8*"
$ws.Range("G6").Value = "0.9"
$ws.Range("G7").Value = "0.5-0.9*"
$ws.Range("G8").Value = "0.2-0.7"
$ws.Range("G9").Value = "N/A"

# --- Column H (new column) ---
$ws.Range("H2").Value = "N/A"
$ws.Range("H3").Value = "N/A"
$ws.Range("H4").Value = "N/A"
$ws.Range("H5").Value = "0.1,0.5"
$ws.Range("H6").Value = "0.7,0.8"
$ws.Range("H7").Value = "0.5-0.9*"
$ws.Range("H8").Value = "0.1-0.4"
$ws.Range("H9").Value = "0.4-0.9*"
